$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invoices")

# Change D14 fill from blue (FF0070C0) to light green (FF92D050), matching B3/B4/B5... style
$ws.Range("D14").Interior.Color = 5296274

# Add yellow fill (FFFFFF00) to B15, matching B13 style
$ws.Range("B15").Interior.Color = 65535

# New content a couple rows below, row 17
$ws.Range("B17").Value = "ADD RETURN"
$ws.Range("C17").Value = "in invoice layout"

# Update selection to B15 to match the new active cell in the saved file
$ws.Range("B15").Select()
